# Apply the cryptos list update (prices / 1h volumes / two name-link row swaps)
# D column holds text-formatted numbers (mirrors the source feed, e.g. "73.587.92"),
# so numeric-looking updates are written through a Text number format and then
# ClearFormats() is used to drop the temporary format, leaving the cell as plain text
# with no style override (matching the original un-styled inline-string cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "73.587.92"; E = "  +6.80%  " },
    @{ Row = 3; D = "2.603.07"; E = "  +6.68%  " },
    @{ Row = 4; D = "0.999"; E = "  -0.08%  " },
    @{ Row = 5; B = "BNB"; C = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"; D = "585.85"; E = "  +4.53%  " },
    @{ Row = 6; B = "Solana"; C = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"; D = "184.29"; E = "  +13.67%  " },
    @{ Row = 7; D = "0.999"; E = "  -0.12%  " },
    @{ Row = 8; D = "0.537"; E = "  +4.28%  " },
    @{ Row = 9; D = "0.201"; E = "  +19.29%  " },
    @{ Row = 10; D = "2.596.57"; E = "  +6.41%  " },
    @{ Row = 11; E = "  +0.11%  " },
    @{ Row = 12; E = "  +9.09%  " },
    @{ Row = 13; E = "  +3.87%  " },
    @{ Row = 14; D = "3.112.61"; E = "  +7.85%  " },
    @{ Row = 15; B = "ShibaInu"; C = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"; D = "0.0000188"; E = "  +6.88%  " },
    @{ Row = 16; B = "WrappedBTC"; C = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"; D = "73.623.55"; E = "  +7.04%  " },
    @{ Row = 17; D = "26.10"; E = "  +12.61%  " },
    @{ Row = 18; D = "2.621.91"; E = "  +7.87%  " },
    @{ Row = 19; D = "8.85"; E = "  +27.48%  " },
    @{ Row = 20; D = "11.82"; E = "  +12.47%  " },
    @{ Row = 21; D = "373.29"; E = "  +10.23%  " },
    @{ Row = 22; E = "  +16.68%  " },
    @{ Row = 23; E = "  +6.52%  " },
    @{ Row = 24; E = "  -0.09%  " },
    @{ Row = 25; D = "69.90"; E = "  +4.25%  " },
    @{ Row = 26; E = "  +12.64%  " },
    @{ Row = 27; D = "9.27"; E = "  +12.81%  " },
    @{ Row = 28; D = "2.733.16"; E = "  +7.02%  " },
    @{ Row = 29; D = "0.998"; E = "  -0.68%  " },
    @{ Row = 30; E = "  +15.24%  " },
    @{ Row = 31; B = "InternetComputer(DFINITY)"; C = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; D = "7.94"; E = "  +11.35%  " },
    @{ Row = 32; B = "Fetch.AI"; C = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"; D = "1.38"; E = "  +19.75%  " },
    @{ Row = 33; D = "505.75"; E = "  +18.10%  " },
    @{ Row = 34; D = "1.74"; E = "  +8.14%  " },
    @{ Row = 35; D = "0.998"; E = "  -0.19%  " },
    @{ Row = 36; D = "0.121"; E = "  +14.19%  " },
    @{ Row = 37; D = "159.76"; E = "  +0.14%  " },
    @{ Row = 38; D = "19.18"; E = "  +6.77%  " },
    @{ Row = 39; D = "19.30"; E = "  +1.41%  " },
    @{ Row = 40; E = "  +0.02%  " },
    @{ Row = 41; E = "  +12.55%  " },
    @{ Row = 42; E = "  +29.98%  " },
    @{ Row = 43; E = "  +10.69%  " },
    @{ Row = 44; D = "0.324"; E = "  +8.98%  " },
    @{ Row = 45; D = "157.30"; E = "  +20.46%  " },
    @{ Row = 46; E = "  +16.68%  " },
    @{ Row = 47; D = "1.17"; E = "  +8.76%  " },
    @{ Row = 48; D = "38.69"; E = "  +3.30%  " },
    @{ Row = 49; D = "3.62"; E = "  +8.18%  " },
    @{ Row = 50; D = "0.525"; E = "  +9.15%  " },
    @{ Row = 51; D = "20.33"; E = "  +20.37%  " }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("B")) { $ws.Cells.Item($u.Row, 2).Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Cells.Item($u.Row, 3).Value = $u.C }
    if ($u.ContainsKey("D")) {
        $cell = $ws.Cells.Item($u.Row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.ClearFormats()
    }
    if ($u.ContainsKey("E")) { $ws.Cells.Item($u.Row, 5).Value = $u.E }
}
